# Auto-generated edit script: updates market-price / profit figures
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets,
# matching a scheduled-runner refresh of current market board prices.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")

$values = @{
    "I6" = 128212.5
    "L6" = 1649.625
    "M6" = -384525.5
    "J6" = 549.875
    "N6" = -1873.625
    "K6" = 384637.5
    "H6" = 64381.188
    "H8" = 1725.6842
    "K8" = 108
    "M8" = 31
    "I8" = 36
    "J87" = 0
    "H87" = 0
    "H90" = 0
    "J90" = 0
    "M98" = 443.0667000000001
    "K98" = 1054.9333
    "I98" = 1054.9333
    "N98" = -3996
    "J98" = 1000
    "H98" = 1051.5
    "L98" = 1000
    "M122" = -714.7999
    "L122" = 3000
    "J122" = 1000
    "I122" = 1054.9333
    "K122" = 3164.7999
    "N122" = -7900
    "H122" = 1051.5
    "J127" = 2996.6667
    "L127" = 8990.000100000001
    "N127" = -18910.0001
    "H127" = 1224942.5
    "N129" = -5015357.800000001
    "M129" = -3748636
    "K129" = 3753636
    "L129" = 5005357.800000001
    "H129" = 1430029.4
    "I129" = 1251212
    "J129" = 1668452.6
    "L138" = 17123.4
    "N138" = -27403.4
    "J138" = 5707.8
    "H138" = 4912.4375
    "L141" = 4563
    "H141" = 1756.72
    "M141" = -266.9500000000007
    "N141" = -14923
    "I141" = 1815.65
    "K141" = 5446.950000000001
    "J141" = 1521
    "L87" = 0
    "L90" = 0
}
foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
$ws.Range("N87").ClearContents()
$ws.Range("N90").ClearContents()

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")

$values = @{
    "I61" = 3807.2
    "K61" = 3807.2
    "H61" = 12860719
    "M61" = -3595.2
    "M97" = -83.59375
    "I97" = 579.59375
    "K97" = 579.59375
    "H97" = 558.4706
    "I102" = 828.6539
    "K102" = 828.6539
    "M102" = 793.3461
    "H102" = 823.8889
    "J135" = 84110.234
    "N135" = -94250.234
    "L135" = 84110.234
    "H135" = 84110.234
    "M136" = -8871.599999999999
    "I136" = 3807.2
    "K136" = 11421.6
    "H136" = 12860719
}
foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")

$values = @{
    "L134" = 18768.8568
    "H134" = 2885.081
    "J134" = 6256.2856
    "N134" = -23838.8568
    "M134" = -3760.399800000001
    "I134" = 2098.4666
    "K134" = 6295.399800000001
}
foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")

$values = @{
    "J58" = 3480.762
    "L58" = 3480.762
    "M58" = -3023.8
    "K58" = 3226.8
    "I58" = 3226.8
    "H58" = 3398.8386
    "N58" = -3886.762
    "J62" = 4886.3335
    "L62" = 4886.3335
    "M62" = -3991.4
    "H62" = 4763.1816
    "N62" = -6134.3335
    "I62" = 4615.4
    "K62" = 4615.4
    "L65" = 24431.6675
    "M65" = -19957
    "K65" = 23077
    "N65" = -30671.6675
    "H65" = 4763.1816
    "I65" = 4615.4
    "J65" = 4886.3335
    "M136" = -7130.400000000001
    "K136" = 9680.400000000001
    "L136" = 10442.286
    "I136" = 3226.8
    "H136" = 3398.8386
    "J136" = 3480.762
    "N136" = -15542.286
}
foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")

$values = @{
    "L113" = 2613.8571
    "H113" = 778.4
    "N113" = -6953.8571
    "J113" = 871.2857
    "K130" = 9000
    "I130" = 3000
    "H130" = 3000
    "M130" = -3980
}
foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")

$values = @{
    "M80" = -929.75
    "K80" = 1927.75
    "H80" = 2649.0557
    "J80" = 4091.6667
    "N80" = -6087.6667
    "L80" = 4091.6667
    "I80" = 1927.75
    "J83" = 4091.6667
    "M83" = -4646.75
    "I83" = 1927.75
    "H83" = 2649.0557
    "K83" = 9638.75
    "N83" = -30442.3335
    "L83" = 20458.3335
    "M122" = -207412
    "L122" = 5391
    "J122" = 1797
    "I122" = 69954
    "K122" = 209862
    "N122" = -10291
    "H122" = 65694.19
    "J123" = 0
    "H123" = 0
    "L123" = 0
}
foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
$ws.Range("N123").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")

$values = @{
    "I7" = 3249.375
    "K7" = 3249.375
    "H7" = 3798.7856
    "M7" = -3137.375
    "J22" = 3569
    "L22" = 3569
    "N22" = -4159
    "H22" = 2711.361
    "H27" = 2711.361
    "L27" = 3569
    "N27" = -3783
    "J27" = 3569
    "K40" = 3988
    "M40" = -3852
    "H40" = 4739.75
    "I40" = 3988
    "H55" = 1121
    "M55" = -9.099999999999994
    "K55" = 182.1
    "I55" = 182.1
    "K126" = 9748.125
    "H126" = 3798.7856
    "I126" = 3249.375
    "M126" = -7278.125
    "H132" = 8344750
    "K132" = 50053500
    "M132" = -50050970
    "I132" = 16684500
    "M136" = -17652295.5
    "I136" = 5884948.5
    "K136" = 17654845.5
    "H136" = 3453218.5
}
foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")

$values = @{
    "L136" = 10982.0001
    "N136" = -16082.0001
    "H136" = 66669944
    "J136" = 3660.6667
}
foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
